# correcao do faturamento diario por lojas
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Bibi Cell Mundi (row 2)
$ws.Range("AE2").Value = 16623.32
$ws.Range("AG2").Value = 301430.49

# Bibi Cell Vieiralves (row 3)
$ws.Range("AE3").Value = 4728.9
$ws.Range("AG3").Value = 154689.89

# Bibi Cell Ponta Negra (row 4)
$ws.Range("AE4").Value = 1399
$ws.Range("AG4").Value = 89679.71000000001

# Bibi Cell Manauara (row 5)
$ws.Range("AE5").Value = 1577.9
$ws.Range("AG5").Value = 88496.14999999999

# total (row 6)
$ws.Range("AE6").Value = 24329.12
$ws.Range("AG6").Value = 634296.24
